$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching original inline-string formatting)
# by forcing Text number format prior to assignment, so Excel doesn't coerce them to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.946.33"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.196.54"
$ws.Range("E3").Value = "  -2.44%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "295.23"
$ws.Range("E5").Value = "  -4.15%  "
$ws.Range("D6").Value = "88.98"
$ws.Range("E6").Value = "  -6.25%  "
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -8.55%  "
$ws.Range("D10").Value = "32.10"
$ws.Range("E10").Value = "  -8.41%  "
$ws.Range("D11").Value = "0.0772"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "6.78"
$ws.Range("E13").Value = "  -6.24%  "
$ws.Range("D14").Value = "2.533.14"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "2.273.78"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").Value = "13.09"
$ws.Range("E16").Value = "  -4.72%  "
$ws.Range("E17").Value = "  -8.28%  "
$ws.Range("D18").Value = "43.630.51"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "0.0₃0887"
$ws.Range("E19").Value = "  -8.35%  "
$ws.Range("D20").Value = "5.84"
$ws.Range("E20").Value = "  -8.96%  "
$ws.Range("D21").Value = "10.81"
$ws.Range("E21").Value = "  -13.54%  "
$ws.Range("D22").Value = "63.03"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("D23").Value = "232.75"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("E24").Value = "  -8.70%  "
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "1.83"
$ws.Range("E26").Value = "  -9.18%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "36.15"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("D29").Value = "9.24"
$ws.Range("E29").Value = "  -6.65%  "
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("D31").Value = "148.46"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  -11.67%  "
$ws.Range("E33").Value = "  -5.36%  "
$ws.Range("D34").Value = "0.0735"
$ws.Range("E34").Value = "  -8.64%  "
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").Value = "2.83"
$ws.Range("E36").Value = "  -9.42%  "
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  -10.09%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0282"
$ws.Range("E39").Value = "  -7.55%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  -7.93%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "3.08"
$ws.Range("E41").Value = "  -11.78%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.01"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "13.00"
$ws.Range("E43").Value = "  -13.16%  "
$ws.Range("D44").Value = "1.785.36"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  +11.68%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "72.72"
$ws.Range("E47").Value = "  -10.23%  "
$ws.Range("D48").Value = "0.173"
$ws.Range("E48").Value = "  -10.84%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "91.91"
$ws.Range("E49").Value = "  -8.29%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "4.53"
$ws.Range("E50").Value = "  -8.82%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "64.75"
$ws.Range("E51").Value = "  -9.04%  "
